$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Add a new row (22) documenting "84. Largest Rectangle in Histogram"
# under the "Stack" category, matching the style pattern already used
# by the adjacent rows (A: plain/bordered, B: red highlight, C: bold).
# ---------------------------------------------------------------------

# Copy base cell formatting (border/wrap/alignment) from row 20 onto row 22
$ws.Range("A20:C20").Copy() | Out-Null
$ws.Range("A22:C22").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Column B uses the red-fill "title" style (same as B15)
$ws.Range("B15").Copy() | Out-Null
$ws.Range("B22").PasteSpecial(-4122) | Out-Null       # xlPasteFormats

# Column C uses the bold header style (same as A1)
$ws.Range("A1").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null       # xlPasteFormats

$excel.CutCopyMode = 0

# Cell values
$ws.Range("A22").Value = "Stack"
$ws.Range("B22").Value = "84. Largest Rectangle in Histogram"

$seg1 = "U need to understand monotonic stack first!!!,,,, "
$seg2 = "This solution is too hard for me to explain watch this video for visual explaination:"
$seg3 = "`n"
$seg4 = "https://youtu.be/zx5Sw9130L0?si=79lWeV1xHcdTWbPt"
$fullText = $seg1 + $seg2 + $seg3 + $seg4
$ws.Range("C22").Value = $fullText

# Rich-text run formatting for C22 (seg1 keeps the default/unformatted
# run; seg2 and seg4 are regular; seg3, the line break, is bold)
$run2 = $ws.Range("C22").Characters(51, 85)
$run2.Font.Name = "Calibri"
$run2.Font.Size = 11
$run2.Font.Bold = $false

$run3 = $ws.Range("C22").Characters(136, 1)
$run3.Font.Name = "Calibri"
$run3.Font.Size = 11
$run3.Font.Bold = $true

$run4 = $ws.Range("C22").Characters(137, 48)
$run4.Font.Name = "Calibri"
$run4.Font.Size = 11
$run4.Font.Bold = $false

# Match the row height used by the source workbook for this entry
$ws.Rows.Item(22).RowHeight = 57.6

# Update the view: scroll so row 11 is at the top and select C23
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C23").Select() | Out-Null

Write-Host "Added row 22: 84. Largest Rectangle in Histogram"
